# Auto-generated Excel COM-interop script to apply cryptos.xlsx price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.825.63"
$ws.Range("E2").Value = "  +0.43%  "
$ws.Range("D3").Value = "1.650.41"
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("E4").Value = "  +0.63%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.506"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.42%  "
$ws.Range("E7").Value = "  +0.57%  "
$ws.Range("E8").Value = "  +0.59%  "
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("E10").Value = "  +0.47%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0844"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.12%  "
$ws.Range("D12").Value = "1.874.77"
$ws.Range("E12").Value = "  -0.19%  "
$ws.Range("D13").Value = "1.656.43"
$ws.Range("E13").Value = "  +0.49%  "
$ws.Range("E14").Value = "  +1.75%  "
$ws.Range("E15").Value = "  +0.44%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.95"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.39%  "
$ws.Range("D17").Value = "26.818.61"
$ws.Range("E17").Value = "  +0.11%  "
$ws.Range("D18").Value = "0.0₃0747"
$ws.Range("E18").Value = "  -0.18%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "217.35"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.53%  "
$ws.Range("E20").Value = "  +0.67%  "
$ws.Range("E21").Value = "  +0.48%  "
$ws.Range("E22").Value = "  +14.25%  "
$ws.Range("E23").Value = "  -0.79%  "
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.96"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.31%  "
$ws.Range("E26").Value = "  +0.64%  "
$ws.Range("E27").Value = "  -0.39%  "
$ws.Range("E28").Value = "  +3.83%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.84"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.27%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0523"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.10%  "
$ws.Range("E31").Value = "  +0.81%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.36"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.57%  "
$ws.Range("E33").Value = "  +0.97%  "
$ws.Range("B34").Value = "Maker"
$ws.Range("C34").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D34").Value = "1.280.89"
$ws.Range("E34").Value = "  -0.18%  "
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.55"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.91%  "
$ws.Range("E36").Value = "  +1.98%  "
$ws.Range("E37").Value = "  -0.54%  "
$ws.Range("E38").Value = "  +5.19%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.837"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.01%  "
$ws.Range("E41").Value = "  +1.73%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.25"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.58%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.44"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.90%  "
$ws.Range("D44").Value = "1.799.95"
$ws.Range("E44").Value = "  +0.64%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.22"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.37%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "59.71"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.50%  "
$ws.Range("E47").Value = "  +1.37%  "
$ws.Range("E48").Value = "  +1.52%  "
$ws.Range("E49").Value = "  +0.50%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.83"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.55%  "
$ws.Range("E51").Value = "  +1.77%  "
